$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data now has two sending/target clusters ("MuSCs" and
# "Resolving-Mac") instead of just one, so the Epcam-Epcam sheet grows
# from a single data row to the full 2x2 cross product of clusters
# (rows 2-5). Columns stay A:T:
#   A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
#   E..T the various expression / specificity metrics.
$data = @(
    @("MuSCs","Epcam","Epcam","MuSCs",1,0.5,0.050085,0.10017,0.4119961612284069,0.318389142285015,1,0.5,0.050085,0.10017,0.4119961612284069,0.318389142285015,0.002508507225,0.0100340289,0.1697408368669435,0.1013716459249875),
    @("MuSCs","Epcam","Epcam","Resolving-Mac",1,0.5,0.050085,0.10017,0.4119961612284069,0.318389142285015,1,0.3333333333333333,0.07148166666666667,0.214445,0.5880038387715931,0.681610857714985,0.003580159275,0.02148095565,0.2422553243614635,0.2170174963600275),
    @("Resolving-Mac","Epcam","Epcam","MuSCs",1,0.3333333333333333,0.07148166666666667,0.214445,0.5880038387715931,0.681610857714985,1,0.5,0.050085,0.10017,0.4119961612284069,0.318389142285015,0.003580159275,0.02148095565,0.2422553243614635,0.2170174963600275),
    @("Resolving-Mac","Epcam","Epcam","Resolving-Mac",1,0.3333333333333333,0.07148166666666667,0.214445,0.5880038387715931,0.681610857714985,1,0.3333333333333333,0.07148166666666667,0.214445,0.5880038387715931,0.681610857714985,0.005109628669444444,0.045986658025,0.3457485144101297,0.4645933613549575)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $i + 2
    $rowValues = $data[$i]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($rowIndex, $c + 1).Value = $rowValues[$c]
    }
}
